# Add a new "2023" column (S) to the trade worksheet, mirroring the
# existing year columns (B:R = 2006:2022).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column R into column S so the styles (s="5"/"6"/"13"
# etc.) carry over to the new column. Done in two pieces (row 1, then rows
# 3:14) so the untouched gap at row 2 is not materialized as a new row.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122) | Out-Null

$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header cell S1 is part of the merged title band A1:R1 -> A1:S1
$ws.Range("A1:S1").Merge()

# Year header row
$ws.Range("S3").Value = 2023

# Data rows 4-14
$ws.Range("S4").Value = 810.5
$ws.Range("S5").Value = 135.19999999999999
$ws.Range("S6").Value = 3146
$ws.Range("S7").Value = 2339
$ws.Range("S8").Value = 1503.1
$ws.Range("S9").Value = 41.9
$ws.Range("S10").Value = 42.8
$ws.Range("S11").Value = 93.3
$ws.Range("S12").Value = 5.9
$ws.Range("S13").Value = 746
$ws.Range("S14").Value = 704.7

# Restore the selection/zoom as Excel would leave them after widening the
# selection to include the new column.
$ws.Range("A1:S1").Select()
$ws.Application.ActiveWindow.Zoom = 100
